# Commit de atualização e correção 24/10/23
#
# The "id" (column A) recorded for one observation (row 30, the second
# tempo/bebida reading of 3.28 / alcool) was mistyped as 9 - a duplicate
# of the id already used in row 17 - when it should have been 29. Fix the
# typo, then (as the author did through Excel's Data > Sort) sort the
# A2:C31 data range by column A ascending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the mistyped id value in A30 (was 9, should be 29) ---
$ws.Range("A30").Value2 = 29

# --- Select the id column, like a user would before sorting ---
[void]$ws.Range("A2:A31").Select()

# --- Sort A2:C31 by column A (id), ascending, no header row in range ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
[void]$sortObj.SortFields.Add($ws.Range("A2:A31"))
$sortObj.SetRange($ws.Range("A2:C31"))
$sortObj.Header = 2
$sortObj.Orientation = 1
$sortObj.Apply()
